$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8831.833000000001
$ws.Range("I62").Value = 6000
$ws.Range("J62").Value = 9398.200000000001
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 9398.200000000001
$ws.Range("M62").Value = -5376
$ws.Range("N62").Value = -10646.2
$ws.Range("H65").Value = 8831.833000000001
$ws.Range("I65").Value = 6000
$ws.Range("J65").Value = 9398.200000000001
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 46991
$ws.Range("M65").Value = -26880
$ws.Range("N65").Value = -53231
$ws.Range("H70").Value = 1549.6666
$ws.Range("J70").Value = 1666
$ws.Range("L70").Value = 4998
$ws.Range("N70").Value = -5538
$ws.Range("H73").Value = 1549.6666
$ws.Range("J73").Value = 1666
$ws.Range("L73").Value = 4998
$ws.Range("N73").Value = -6870
$ws.Range("H74").Value = 25004780
$ws.Range("J74").Value = 25004780
$ws.Range("L74").Value = 25004780
$ws.Range("N74").Value = -25006652
$ws.Range("H76").Value = 1987265
$ws.Range("I76").Value = 3207.8572
$ws.Range("J76").Value = 7939436.5
$ws.Range("K76").Value = 3207.8572
$ws.Range("L76").Value = 7939436.5
$ws.Range("M76").Value = -2892.8572
$ws.Range("N76").Value = -7940066.5
$ws.Range("H77").Value = 25004780
$ws.Range("J77").Value = 25004780
$ws.Range("L77").Value = 125023900
$ws.Range("N77").Value = -125033260
$ws.Range("H79").Value = 1987265
$ws.Range("I79").Value = 3207.8572
$ws.Range("J79").Value = 7939436.5
$ws.Range("K79").Value = 3207.8572
$ws.Range("L79").Value = 7939436.5
$ws.Range("M79").Value = -2115.8572
$ws.Range("N79").Value = -7941620.5
$ws.Range("H86").Value = 12204.556
$ws.Range("I86").Value = 1183.1666
$ws.Range("J86").Value = 34247.332
$ws.Range("K86").Value = 1183.1666
$ws.Range("L86").Value = 34247.332
$ws.Range("M86").Value = -60.16660000000002
$ws.Range("N86").Value = -36493.332
$ws.Range("H89").Value = 12204.556
$ws.Range("I89").Value = 1183.1666
$ws.Range("J89").Value = 34247.332
$ws.Range("K89").Value = 5915.833000000001
$ws.Range("L89").Value = 171236.66
$ws.Range("M89").Value = -299.8330000000005
$ws.Range("N89").Value = -182468.66
$ws.Range("H99").Value = 209.83333
$ws.Range("J99").Value = 259
$ws.Range("L99").Value = 777
$ws.Range("N99").Value = -3773
$ws.Range("H101").Value = 294.44446
$ws.Range("I101").Value = 258.33334
$ws.Range("J101").Value = 366.66666
$ws.Range("K101").Value = 775.0000200000001
$ws.Range("L101").Value = 1099.99998
$ws.Range("M101").Value = 846.9999799999999
$ws.Range("N101").Value = -4343.999980000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2976.9375
$ws.Range("I45").Value = 2632.111
$ws.Range("J45").Value = 3420.2856
$ws.Range("K45").Value = 2632.111
$ws.Range("L45").Value = 3420.2856
$ws.Range("M45").Value = -2255.111
$ws.Range("N45").Value = -4174.2856
$ws.Range("H97").Value = 250002750
$ws.Range("I97").Value = 5255
$ws.Range("J97").Value = 500000260
$ws.Range("K97").Value = 5255
$ws.Range("L97").Value = 500000260
$ws.Range("M97").Value = -4759
$ws.Range("N97").Value = -500001252
$ws.Range("H102").Value = 979.44446
$ws.Range("I102").Value = 990
$ws.Range("J102").Value = 926.6667
$ws.Range("K102").Value = 990
$ws.Range("L102").Value = 926.6667
$ws.Range("M102").Value = 632
$ws.Range("N102").Value = -4170.6667
$ws.Range("H122").Value = 4091.5
$ws.Range("I122").Value = 3157.1428
$ws.Range("J122").Value = 5399.6
$ws.Range("K122").Value = 9471.428400000001
$ws.Range("L122").Value = 16198.8
$ws.Range("M122").Value = -7021.428400000001
$ws.Range("N122").Value = -21098.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1103.5385
$ws.Range("I94").Value = 893.25
$ws.Range("J94").Value = 1440
$ws.Range("K94").Value = 893.25
$ws.Range("L94").Value = 1440
$ws.Range("M94").Value = -442.25
$ws.Range("N94").Value = -2342
$ws.Range("H99").Value = 1352.0667
$ws.Range("J99").Value = 1302.2
$ws.Range("L99").Value = 1302.2
$ws.Range("N99").Value = -4298.2
$ws.Range("H105").Value = 4169110
$ws.Range("I105").Value = 2317.1428
$ws.Range("J105").Value = 10002620
$ws.Range("K105").Value = 2317.1428
$ws.Range("L105").Value = 10002620
$ws.Range("M105").Value = -570.1428000000001
$ws.Range("N105").Value = -10006114
$ws.Range("H107").Value = 788.75
$ws.Range("I107").Value = 810
$ws.Range("J107").Value = 753.3333
$ws.Range("K107").Value = 810
$ws.Range("L107").Value = 753.3333
$ws.Range("M107").Value = 1110
$ws.Range("N107").Value = -4593.3333
$ws.Range("H134").Value = 6812.7144
$ws.Range("I134").Value = 6812.7144
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 20438.1432
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -17903.1432
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1700
$ws.Range("I134").Value = 1600
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 4800
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -2265
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 34167.668
$ws.Range("I68").Value = 1500
$ws.Range("K68").Value = 4500
$ws.Range("M68").Value = -3689
$ws.Range("H71").Value = 34167.668
$ws.Range("I71").Value = 1500
$ws.Range("K71").Value = 13500
$ws.Range("M71").Value = -9444
$ws.Range("H131").Value = 710.09
$ws.Range("J131").Value = 721.46313
$ws.Range("L131").Value = 2164.38939
$ws.Range("N131").Value = -12244.38939
$ws.Range("H132").Value = 593.6
$ws.Range("I132").Value = 617
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 5553
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -3023
$ws.Range("N132").Value = -9560
$ws.Range("H140").Value = 2506.087
$ws.Range("I140").Value = 1235.3846
$ws.Range("K140").Value = 3706.1538
$ws.Range("M140").Value = 1473.8462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 33337058
$ws.Range("I102").Value = 45458404
$ws.Range("J102").Value = 3353.5
$ws.Range("K102").Value = 45458404
$ws.Range("L102").Value = 3353.5
$ws.Range("M102").Value = -45456782
$ws.Range("N102").Value = -6597.5
$ws.Range("H113").Value = 2084.625
$ws.Range("I113").Value = 1673
$ws.Range("J113").Value = 2907.875
$ws.Range("K113").Value = 1673
$ws.Range("L113").Value = 2907.875
$ws.Range("M113").Value = 497
$ws.Range("N113").Value = -7247.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1056.5714
$ws.Range("I46").Value = 874.5
$ws.Range("J46").Value = 1299.3334
$ws.Range("K46").Value = 874.5
$ws.Range("L46").Value = 1299.3334
$ws.Range("M46").Value = -686.5
$ws.Range("N46").Value = -1675.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2815.25
$ws.Range("I126").Value = 2816.6667
$ws.Range("J126").Value = 2802.5
$ws.Range("K126").Value = 8450.000100000001
$ws.Range("L126").Value = 8407.5
$ws.Range("M126").Value = -5980.000100000001
$ws.Range("N126").Value = -13347.5
